$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 holds the single candidate record; update the generated
# credentials/id fields per the new iProctor registration batch.
$ws.Range("A2").Value = "OOkLX161"
$ws.Range("B2").Value = 23071823
$ws.Range("C2").Value = "urfnboj78"
$ws.Range("D2").Value = "KN&8w7j%"
$ws.Range("F2").Value = "RXpSxSmD"
$ws.Range("G2").Value = "reGn"
